# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on the
# per-language sheets, as part of regenerating the handback status
# report.
#
# Only row 2 (the 1dc12252-...-82d5-2b0548f8ed8d handback entry) changes
# on both the "zh-cn" and "de-de" sheets:
#   zh-cn: E2 14:50:35 -> 14:51:16 ; H2 14:50:51 -> 14:51:34
#   de-de: E2 14:50:38 -> 14:51:19 ; H2 14:50:58 -> 14:51:39

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-12 14:51:16"
$wsZh.Range("H2").Value = "2016-03-12 14:51:34"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-12 14:51:19"
$wsDe.Range("H2").Value = "2016-03-12 14:51:39"
